# Update the "想去人数" (interested-people count) column F for a handful of
# events on the 展览 (Exhibition) and 全部类型 (All Types) sheets.
# Same underlying rows are duplicated across both sheets, with 全部类型
# having one extra row (the 演出 entry) inserted before the last row.

$wb = $excel.ActiveWorkbook

$updates = @{
    "F2"  = 314
    "F3"  = 56
    "F4"  = 477
    "F5"  = 4590
    "F6"  = 355
    "F7"  = 625
    "F9"  = 709
    "F10" = 197
}

$ws1 = $wb.Worksheets.Item("展览")
foreach ($addr in $updates.Keys) {
    $ws1.Range($addr).Value = $updates[$addr]
}

$updates4 = @{
    "F2"  = 314
    "F3"  = 56
    "F4"  = 477
    "F5"  = 4590
    "F6"  = 355
    "F7"  = 625
    "F9"  = 709
    "F11" = 197
}

$ws4 = $wb.Worksheets.Item("全部类型")
foreach ($addr in $updates4.Keys) {
    $ws4.Range($addr).Value = $updates4[$addr]
}
